$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7115757465362549
$ws.Range("B1").Value = 1.462210178375244
$ws.Range("C1").Value = 4.195500373840332
$ws.Range("D1").Value = 2.405681133270264
$ws.Range("E1").Value = 0.5585829615592957
